{"js": "const replacements = [\n  [\"2025-11-26 Wednesday\", \"2025-11-27 Thursday\"],\n  [\"334\u00d72=\", \"530\u00d76=\"],\n  [\"563\u00d79=\", \"853\u00d77=\"],\n  [\"566\u00d74=\", \"373\u00d79=\"],\n  [\"947\u00d79=\", \"563\u00d72=\"],\n  [\"576\u00d77=\", \"636\u00d73=\"],\n  [\"481\u00d72=\", \"872\u00d77=\"],\n  [\"389\u00d72=\", \"477\u00d79=\"],\n  [\"405\u00d78=\", \"902\u00d72=\"],\n  [\"272\u00d72=\", \"221\u00d75=\"],\n  [\"820\u00d75=\", \"864\u00d77=\"],\n  [\"231\u00d77=\", \"794\u00d75=\"],\n  [\"298\u00d77=\", \"635\u00d74=\"],\n  [\"135\u00d77=\", \"718\u00d74=\"],\n  [\"973\u00d72=\", \"251\u00d79=\"],\n  [\"950\u00d77=\", \"730\u00d77=\"],\n  [\"709\u00d78=\", \"239\u00d76=\"],\n  [\"846\u00d78=\", \"789\u00d73=\"],\n  [\"274\u00d75=\", \"763\u00d75=\"],\n  [\"259\u00d72=\", \"617\u00d75=\"],\n  [\"957\u00d77=\", \"181\u00d74=\"],\n  [\"347\u00d73=\", \"913\u00d77=\"],\n  [\"911\u00d72=\", \"417\u00d75=\"],\n  [\"788\u00d76=\", \"737\u00d75=\"],\n  [\"707\u00d77=\", \"337\u00d72=\"],\n  [\"715\u00d79=\", \"495\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old='2025-11-26 Wednesday'; New='2025-11-27 Thursday'},\n    @{Old='334\u00d72='; New='530\u00d76='},\n    @{Old='563\u00d79='; New='853\u00d77='},\n    @{Old='566\u00d74='; New='373\u00d79='},\n    @{Old='947\u00d79='; New='563\u00d72='},\n    @{Old='576\u00d77='; New='636\u00d73='},\n    @{Old='481\u00d72='; New='872\u00d77='},\n    @{Old='389\u00d72='; New='477\u00d79='},\n    @{Old='405\u00d78='; New='902\u00d72='},\n    @{Old='272\u00d72='; New='221\u00d75='},\n    @{Old='820\u00d75='; New='864\u00d77='},\n    @{Old='231\u00d77='; New='794\u00d75='},\n    @{Old='298\u00d77='; New='635\u00d74='},\n    @{Old='135\u00d77='; New='718\u00d74='},\n    @{Old='973\u00d72='; New='251\u00d79='},\n    @{Old='950\u00d77='; New='730\u00d77='},\n    @{Old='709\u00d78='; New='239\u00d76='},\n    @{Old='846\u00d78='; New='789\u00d73='},\n    @{Old='274\u00d75='; New='763\u00d75='},\n    @{Old='259\u00d72='; New='617\u00d75='},\n    @{Old='957\u00d77='; New='181\u00d74='},\n    @{Old='347\u00d73='; New='913\u00d77='},\n    @{Old='911\u00d72='; New='417\u00d75='},\n    @{Old='788\u00d76='; New='737\u00d75='},\n    @{Old='707\u00d77='; New='337\u00d72='},\n    @{Old='715\u00d79='; New='495\u00d78='}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
